$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

function Set-PlainCell($ref, $val) {
    $ws.Range($ref).Value = $val
}

# Row 2
Set-TextCell "D2" '65.278.88'
Set-PlainCell "E2" '  +2.38%  '

# Row 3
Set-TextCell "D3" '3.175.05'
Set-PlainCell "E3" '  +3.89%  '

# Row 4
Set-TextCell "D4" '0.999'
Set-PlainCell "E4" '  -0.13%  '

# Row 5
Set-TextCell "D5" '574.62'
Set-PlainCell "E5" '  +3.25%  '

# Row 6
Set-TextCell "D6" '150.97'
Set-PlainCell "E6" '  +6.11%  '

# Row 7
Set-TextCell "D7" '0.999'
Set-PlainCell "E7" '  -0.07%  '

# Row 8
Set-TextCell "D8" '3.169.52'
Set-PlainCell "E8" '  +3.82%  '

# Row 9
Set-TextCell "D9" '0.528'
Set-PlainCell "E9" '  +4.50%  '

# Row 10
Set-TextCell "D10" '0.163'
Set-PlainCell "E10" '  +4.98%  '

# Row 11
Set-TextCell "D11" '6.28'
Set-PlainCell "E11" '  +3.02%  '

# Row 12
Set-TextCell "D12" '0.511'
Set-PlainCell "E12" '  +7.09%  '

# Row 13
Set-TextCell "D13" '0.0000271'
Set-PlainCell "E13" '  +17.20%  '

# Row 14
Set-TextCell "D14" '38.44'
Set-PlainCell "E14" '  +9.89%  '

# Row 15
Set-TextCell "D15" '3.680.18'
Set-PlainCell "E15" '  +3.66%  '

# Row 16
Set-TextCell "D16" '65.251.44'
Set-PlainCell "E16" '  +2.30%  '

# Row 17
Set-TextCell "D17" '7.26'
Set-PlainCell "E17" '  +7.76%  '

# Row 18
Set-TextCell "D18" '3.159.05'
Set-PlainCell "E18" '  +3.54%  '

# Row 19
Set-PlainCell "E19" '  +0.96%  '

# Row 20
Set-TextCell "D20" '515.13'
Set-PlainCell "E20" '  +7.64%  '

# Row 21
Set-TextCell "D21" '15.03'
Set-PlainCell "E21" '  +7.07%  '

# Row 22
Set-TextCell "D22" '16.23'
Set-PlainCell "E22" '  +12.89%  '

# Row 23
Set-TextCell "D23" '0.744'
Set-PlainCell "E23" '  +9.68%  '

# Row 24
Set-TextCell "D24" '7.90'
Set-PlainCell "E24" '  +4.41%  '

# Row 25
Set-TextCell "D25" '85.15'
Set-PlainCell "E25" '  +4.64%  '

# Row 26
Set-TextCell "D26" '1.00'
Set-PlainCell "E26" '  +0.11%  '

# Row 27
Set-TextCell "D27" '9.19'
Set-PlainCell "E27" '  +15.92%  '

# Row 28
Set-TextCell "D28" '2.92'
Set-PlainCell "E28" '  +4.37%  '

# Row 29
Set-TextCell "D29" '2.21'
Set-PlainCell "E29" '  +8.69%  '

# Row 30
Set-TextCell "D30" '28.18'
Set-PlainCell "E30" '  +7.51%  '

# Row 31
Set-TextCell "D31" '2.79'
Set-PlainCell "E31" '  +14.82%  '

# Row 33
Set-TextCell "D33" '0.998'
Set-PlainCell "E33" '  -0.10%  '

# Row 34
Set-TextCell "D34" '6.32'
Set-PlainCell "E34" '  +11.93%  '

# Row 35
Set-TextCell "D35" '6.69'
Set-PlainCell "E35" '  +7.08%  '

# Row 36
Set-TextCell "D36" '55.87'
Set-PlainCell "E36" '  +1.76%  '

# Row 37
Set-TextCell "D37" '478.70'
Set-PlainCell "E37" '  +7.17%  '

# Row 38
Set-TextCell "D38" '0.0879'
Set-PlainCell "E38" '  +8.31%  '

# Row 39
Set-PlainCell "B39" 'dogwifhat'
Set-PlainCell "C39" 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextCell "D39" '3.09'
Set-PlainCell "E39" '  +8.72%  '

# Row 40
Set-PlainCell "B40" 'VeChain'
Set-PlainCell "C40" 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextCell "D40" '0.0423'
Set-PlainCell "E40" '  +3.16%  '

# Row 41
Set-TextCell "D41" '3.128.03'
Set-PlainCell "E41" '  +5.71%  '

# Row 42
Set-TextCell "D42" '8.70'
Set-PlainCell "E42" '  +5.58%  '

# Row 43
Set-TextCell "D43" '0.121'
Set-PlainCell "E43" '  +7.38%  '

# Row 44
Set-PlainCell "B44" 'Fetch.AI'
Set-PlainCell "C44" 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextCell "D44" '2.53'
Set-PlainCell "E44" '  +17.48%  '

# Row 45
Set-PlainCell "B45" 'TheGraph'
Set-PlainCell "C45" 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextCell "D45" '0.292'
Set-PlainCell "E45" '  +12.01%  '

# Row 46
Set-TextCell "D46" '29.42'
Set-PlainCell "E46" '  +5.68%  '

# Row 47
Set-TextCell "D47" '0.0₃0594'
Set-PlainCell "E47" '  +15.07%  '

# Row 48
Set-PlainCell "E48" '  -0.06%  '

# Row 49
Set-PlainCell "B49" 'Stellar'
Set-PlainCell "C49" 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextCell "D49" '0.116'
Set-PlainCell "E49" '  +3.20%  '

# Row 50
Set-PlainCell "B50" 'ThetaToken'
Set-PlainCell "C50" 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
Set-TextCell "D50" '2.34'
Set-PlainCell "E50" '  +12.19%  '

# Row 51
Set-TextCell "D51" '124.75'
Set-PlainCell "E51" '  +6.75%  '
